$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition / resize the "Chart 1" chart object ---
# Original anchor: from col=9(J) colOff=447675, row=2(3), rowOff=171450
#                   to   col=17(R) colOff=142875, row=17(18), rowOff=57150
# New anchor:      from col=1(B) colOff=381000, row=14(15), rowOff=142875
#                   to   col=8(I) colOff=381000, row=29(30), rowOff=28575
# (moved down/left to make room for newly added R and S measurement pictures)
$co = $ws.ChartObjects().Item(1)
$co.Left   = 88.4375
$co.Top    = 221.25
$co.Width  = 426.357421875
$co.Height = 216.0

# --- Touch the X axis line formatting (Reverse Voltage axis) ---
# Excel normalizes the axis line spPr (adds explicit line weight) once the
# axis format is touched, even though the line itself stays hidden (noFill).
$cht = $co.Chart
$xAxis = $cht.Axes(1)
$xAxis.Format.Line.Weight = 0.75
